# Add a new "Статус" (Status) column to the reviews header row.
#
# The header row currently reads:
#   A=Дата, B=Время, C=Точка, D=Имя, E=Фамилия, F=ID
# A new column is inserted at D so the row becomes:
#   A=Дата, B=Время, C=Точка, D=Статус(new), E=Имя, F=Фамилия, G=ID
# i.e. existing columns D:F (Имя/Фамилия/ID) shift one place right and
# keep their original header formatting; the new D1 cell picks up the
# same header style (bold font, thin box border, centered alignment)
# used by the rest of the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new column before column D - this shifts Имя/Фамилия/ID
# (and any data below them) one column to the right automatically.
[void]$ws.Columns("D:D").Insert()

# New header cell + its text.
$ws.Range("D1").Value = "Статус"

# Match the header styling used by the surrounding header cells
# (bold Calibri font, thin border box, centered horizontal/top vertical
# alignment) so the new header looks consistent with the rest of row 1.
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D1").VerticalAlignment = -4160     # xlTop
$ws.Range("D1").Borders.LineStyle = 1         # xlContinuous (thin box border)

# Restore the view's active cell/selection.
[void]$ws.Range("H18").Select()
